$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The to-do list rows 2-8 got re-ordered (row 8's task rotated to the top,
# rows 2-7 shifted down one). Rows 9-10 are unchanged.
$ws.Range("A2").Value = "Error handling strategy in Model Compiler"
$ws.Range("B2").Value = 7

$ws.Range("A3").Value = "Refactoring - we need consistency across the board - All 3 apps"
$ws.Range("B3").Value = 21

$ws.Range("A4").Value = "Make exporter a GUP.  Build UI & hook data into the max files"
$ws.Range("B4").Value = 14

$ws.Range("A5").Value = "Make the path from Max->Model Viewer seamless"
$ws.Range("B5").Value = 4

$ws.Range("A6").Value = "Add full screen support"
$ws.Range("B6").Value = 3

$ws.Range("A7").Value = "Add camera controls to model viewer"
$ws.Range("B7").Value = 5

$ws.Range("A8").Value = "Textured surfaces"
$ws.Range("B8").Value = 35
